$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 1. Everything that used to be on row 1
# (the original text header) shifts down to row 2, and so on through the
# rest of the table (old row 2 -> row 3, ... old row 25 -> row 26).
$ws.Rows.Item(1).Insert()

# The just-shifted former header row (now row 2) carries the bold/border/
# center style that used to live on row 1. Grab that formatting and copy
# it onto the brand-new row 1 before we overwrite row 2's look, so the new
# numeric header row keeps the original "header" style and row 2 goes back
# to being a normal, unstyled row (matching the shifted-down plain data rows).
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1:L1").ClearContents()
$excel.CutCopyMode = $false

$ws.Range("A2:L2").ClearFormats()

# New row 1: plain numeric column-index header (0-based), 0..11.
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Row 2 keeps the original textual header labels (which were on row 1
# before the insert) ...
$ws.Range("A2").Value = "Lg.,mm"
$ws.Range("B2").Value = "Threading"
$ws.Range("C2").Value = "HeadDia., mm"
$ws.Range("D2").Value = "HeadHt., mm"
$ws.Range("E2").Value = "DriveSize"
$ws.Range("F2").Value = "TensileStrength, psi"
$ws.Range("G2").Value = "Specifications Met"
$ws.Range("H2").Value = "Pkg.Qty."
$ws.Range("J2").Value = "Pkg."

# ... except I2, K2 and L2, which are blank in the new layout.
$ws.Range("I2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()

Write-Output "done"
